$d = $word.ActiveDocument

# Update the date/weekday heading (first paragraph).
$d.Paragraphs.Item(1).Range.Text = "2025-01-28 Tuesday"

# Update the division-problem table cells by explicit (row, column)
# position, since several old/new values collide with each other and a
# blind text-based Find/Replace could double-substitute.
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "25÷8="
$tbl.Cell(1, 2).Range.Text = "65÷9="
$tbl.Cell(1, 3).Range.Text = "34÷6="
$tbl.Cell(1, 4).Range.Text = "91÷3="
$tbl.Cell(1, 5).Range.Text = "62÷7="

$tbl.Cell(5, 1).Range.Text = "46÷3="
$tbl.Cell(5, 2).Range.Text = "68÷5="
$tbl.Cell(5, 3).Range.Text = "81÷5="
$tbl.Cell(5, 4).Range.Text = "67÷7="
$tbl.Cell(5, 5).Range.Text = "59÷8="

$tbl.Cell(9, 1).Range.Text = "53÷8="
$tbl.Cell(9, 2).Range.Text = "95÷6="
$tbl.Cell(9, 3).Range.Text = "35÷9="
$tbl.Cell(9, 4).Range.Text = "85÷4="
$tbl.Cell(9, 5).Range.Text = "59÷7="

$tbl.Cell(13, 1).Range.Text = "85÷4="
$tbl.Cell(13, 2).Range.Text = "27÷3="
$tbl.Cell(13, 3).Range.Text = "61÷8="
$tbl.Cell(13, 4).Range.Text = "74÷2="
$tbl.Cell(13, 5).Range.Text = "33÷7="

$tbl.Cell(17, 1).Range.Text = "33÷2="
$tbl.Cell(17, 2).Range.Text = "98÷7="
$tbl.Cell(17, 3).Range.Text = "41÷2="
$tbl.Cell(17, 4).Range.Text = "62÷5="
$tbl.Cell(17, 5).Range.Text = "51÷5="
